$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I0 and IF (copy the header formatting from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-23: I = 1, J = same value as column H on that row
for ($r = 2; $r -le 23; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 24 is a special case with explicit values
$ws.Cells.Item(24, 9).Value = 4
$ws.Cells.Item(24, 10).Value = 5
